# Auto-generated script applying updated market-price / profit values
# to each profession sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# matching the 'chore: update Sheets via scheduled runner' data refresh.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 110.5
$ws.Range("I2").Value = 110.5
$ws.Range("K2").Value = 110.5
$ws.Range("M2").Value = 2.5
$ws.Range("H43").Value = 4098.8237
$ws.Range("I43").Value = 2999.5
$ws.Range("J43").Value = 5076
$ws.Range("K43").Value = 2999.5
$ws.Range("L43").Value = 5076
$ws.Range("M43").Value = -2930.5
$ws.Range("N43").Value = -5214
$ws.Range("H80").Value = 1586.5555
$ws.Range("J80").Value = 3510
$ws.Range("L80").Value = 10530
$ws.Range("N80").Value = -12526
$ws.Range("H83").Value = 1586.5555
$ws.Range("J83").Value = 3510
$ws.Range("L83").Value = 31590
$ws.Range("N83").Value = -41574
$ws.Range("H86").Value = 6016.375
$ws.Range("I86").Value = 3598.3333
$ws.Range("K86").Value = 3598.3333
$ws.Range("M86").Value = -2475.3333
$ws.Range("H89").Value = 6016.375
$ws.Range("I89").Value = 3598.3333
$ws.Range("K89").Value = 17991.6665
$ws.Range("M89").Value = -12375.6665
$ws.Range("H100").Value = 2257.7144
$ws.Range("I100").Value = 2257.7144
$ws.Range("K100").Value = 2257.7144
$ws.Range("M100").Value = -1716.7144
$ws.Range("H107").Value = 567.4545000000001
$ws.Range("I107").Value = 599.2
$ws.Range("K107").Value = 599.2
$ws.Range("M107").Value = 1320.8
$ws.Range("H113").Value = 8441.615
$ws.Range("I113").Value = 8157.8184
$ws.Range("K113").Value = 8157.8184
$ws.Range("M113").Value = -4903.8184
$ws.Range("H116").Value = 4089.8
$ws.Range("J116").Value = 4483
$ws.Range("L116").Value = 4483
$ws.Range("N116").Value = -11367
$ws.Range("H135").Value = 2197.875
$ws.Range("I135").Value = 518.8
$ws.Range("K135").Value = 4669.2
$ws.Range("M135").Value = -2134.2
$ws.Range("H137").Value = 1456.75
$ws.Range("I137").Value = 1519.75
$ws.Range("J137").Value = 1393.75
$ws.Range("K137").Value = 4559.25
$ws.Range("L137").Value = 4181.25
$ws.Range("M137").Value = -2009.25
$ws.Range("N137").Value = -9281.25
$ws.Range("H138").Value = 3959.7
$ws.Range("I138").Value = 3849.5
$ws.Range("J138").Value = 4125
$ws.Range("K138").Value = 11548.5
$ws.Range("L138").Value = 12375
$ws.Range("M138").Value = -6408.5
$ws.Range("N138").Value = -22655
$ws.Range("H141").Value = 772.8
$ws.Range("I141").Value = 772.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2318.4
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2861.6
$ws.Range("N141").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 11986.667
$ws.Range("I28").Value = 11986.667
$ws.Range("K28").Value = 11986.667
$ws.Range("M28").Value = -11794.667
$ws.Range("H32").Value = 5808.533
$ws.Range("I32").Value = 5525.2856
$ws.Range("K32").Value = 5525.2856
$ws.Range("M32").Value = -5238.2856
$ws.Range("H45").Value = 3087.2104
$ws.Range("I45").Value = 1246.4445
$ws.Range("K45").Value = 1246.4445
$ws.Range("M45").Value = -869.4445000000001
$ws.Range("H74").Value = 998
$ws.Range("I74").Value = 998
$ws.Range("K74").Value = 998
$ws.Range("M74").Value = -124
$ws.Range("H77").Value = 998
$ws.Range("I77").Value = 998
$ws.Range("K77").Value = 4990
$ws.Range("M77").Value = -622
$ws.Range("H99").Value = 11986.667
$ws.Range("I99").Value = 11986.667
$ws.Range("K99").Value = 11986.667
$ws.Range("M99").Value = -8991.666999999999
$ws.Range("H102").Value = 4997.2
$ws.Range("I102").Value = 4997.2
$ws.Range("K102").Value = 4997.2
$ws.Range("M102").Value = -3375.2
$ws.Range("H110").Value = 3035.6667
$ws.Range("I110").Value = 3035.6667
$ws.Range("K110").Value = 3035.6667
$ws.Range("M110").Value = -990.6667000000002
$ws.Range("H132").Value = 1294.8334
$ws.Range("I132").Value = 1294.8334
$ws.Range("K132").Value = 3884.5002
$ws.Range("M132").Value = -1354.5002

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1117.5555
$ws.Range("I94").Value = 632.25
$ws.Range("K94").Value = 632.25
$ws.Range("M94").Value = -181.25
$ws.Range("H105").Value = 2894.7646
$ws.Range("I105").Value = 2826.3125
$ws.Range("K105").Value = 2826.3125
$ws.Range("M105").Value = -1079.3125

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2592.5
$ws.Range("I31").Value = 1984
$ws.Range("K31").Value = 1984
$ws.Range("M31").Value = -1689
$ws.Range("H34").Value = 2592.5
$ws.Range("I34").Value = 1984
$ws.Range("K34").Value = 1984
$ws.Range("M34").Value = -1782
$ws.Range("H39").Value = 4275.5
$ws.Range("I39").Value = 4275.5
$ws.Range("K39").Value = 4275.5
$ws.Range("M39").Value = -3884.5
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H49").Value = 4275.5
$ws.Range("I49").Value = 4275.5
$ws.Range("K49").Value = 4275.5
$ws.Range("M49").Value = -4093.5
$ws.Range("H50").Value = 20000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 20000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 20000
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -21250
$ws.Range("H56").Value = 1700
$ws.Range("I56").Value = 1700
$ws.Range("K56").Value = 1700
$ws.Range("M56").Value = -855
$ws.Range("H58").Value = 3006
$ws.Range("I58").Value = 3006
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3006
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2803
$ws.Range("N58").ClearContents()
$ws.Range("H122").Value = 2136.3635
$ws.Range("I122").Value = 2136.3635
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6409.0905
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3959.0905
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 3006
$ws.Range("I136").Value = 3006
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9018
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6468
$ws.Range("N136").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 296.16666
$ws.Range("J12").Value = 355.2
$ws.Range("L12").Value = 1065.6
$ws.Range("N12").Value = -1411.6
$ws.Range("H41").Value = 107.6
$ws.Range("I41").Value = 107.6
$ws.Range("K41").Value = 322.8
$ws.Range("M41").Value = 15.20000000000005
$ws.Range("H81").Value = 100
$ws.Range("I81").Value = 100
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 300
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 823
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 100
$ws.Range("I84").Value = 100
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 900
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 4716
$ws.Range("N84").ClearContents()
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H131").Value = 894.6818
$ws.Range("I131").Value = 931.2
$ws.Range("K131").Value = 2793.6
$ws.Range("M131").Value = 2246.4

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 14994
$ws.Range("J54").Value = 14994
$ws.Range("L54").Value = 14994
$ws.Range("N54").Value = -15774
$ws.Range("H102").Value = 2480.7
$ws.Range("I102").Value = 2302.4285
$ws.Range("K102").Value = 2302.4285
$ws.Range("M102").Value = -680.4285
$ws.Range("H122").Value = 1999.4
$ws.Range("J122").Value = 1498
$ws.Range("L122").Value = 4494
$ws.Range("N122").Value = -9394
$ws.Range("H132").Value = 2140.6667
$ws.Range("I132").Value = 2140.6667
$ws.Range("K132").Value = 6422.000100000001
$ws.Range("M132").Value = -3892.000100000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1950.5
$ws.Range("I7").Value = 1950.5
$ws.Range("K7").Value = 1950.5
$ws.Range("M7").Value = -1838.5
$ws.Range("H126").Value = 1950.5
$ws.Range("I126").Value = 1950.5
$ws.Range("K126").Value = 5851.5
$ws.Range("M126").Value = -3381.5
$ws.Range("H132").Value = 2035.5
$ws.Range("I132").Value = 2242.6
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 6727.799999999999
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -4197.799999999999
$ws.Range("N132").Value = -8060

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1369
$ws.Range("I122").Value = 1336.25
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4008.75
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1558.75
$ws.Range("N122").Value = -9400
$ws.Range("H126").Value = 3451.8
$ws.Range("I126").Value = 2814.75
$ws.Range("K126").Value = 8444.25
$ws.Range("M126").Value = -5974.25
$ws.Range("H132").Value = 2267
$ws.Range("I132").Value = 1950.5
$ws.Range("K132").Value = 5851.5
$ws.Range("M132").Value = -3321.5

